$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Priority column (E) flips from "low" to "ht" for rows 4-7 on both locale sheets
foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}

# Latest Handoff Datetime (H) updated for zh-cn rows 4-7 (new handoff xliff generated)
foreach ($row in 4..7) {
    $zhcn.Range("H$row").Value = "2016-08-27 04:30:29"
}

# Overview sheet: "Ready for handoff" rows get a refreshed Latest HO Xliff Generate Date
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = "2016-08-27 04:30:34"
}

# de-de's Latest Handoff Datetime for rows 4-7 happens to share the same text value
# as the Overview's Latest HO Xliff Generate Date, so it must track the same refresh.
foreach ($row in 4..7) {
    $dede.Range("H$row").Value = "2016-08-27 04:30:34"
}
